# Auto-generated files on 2025-11-11
# Updates the "HotStock_Top20" style sheet: rows 2-21 of columns A, B, C
# are refreshed with the new hot-stock rankings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2..21, columns A (1), B (2), C (3).
$data = @(
    @("合富中国", "合富中国", "平潭发展"),
    @("中国中免", "中国中免", "合富中国"),
    @("方正电机", "平潭发展", "安泰集团"),
    @("平潭发展", "福龙马", "福龙马"),
    @("天际股份", "天际股份", "摩恩电气"),
    @("摩恩电气", "方正电机", "海马汽车"),
    @("上纬新材", "特变电工", "方正电机"),
    @("福龙马", "海马汽车", "香农芯创"),
    @("海马汽车", "东百集团", "中国中免"),
    @("特变电工", "摩恩电气", "天际股份"),
    @("东百集团", "多氟多", "特变电工"),
    @("香农芯创", "黄河旋风", "大有能源"),
    @("安泰集团", "香农芯创", "孚日股份"),
    @("永太科技", "亿晶光电", "雪人集团"),
    @("黄河旋风", "上纬新材", "三花智控"),
    @("亿晶光电", "永太科技", "东百集团"),
    @("四方达", "安泰集团", "黄河旋风"),
    @("多氟多", "协鑫集成", "亿晶光电"),
    @("雪人集团", "中来股份", "上纬新材"),
    @("大有能源", "雪人集团", "海陆重工")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rowValues[0]
    $ws.Cells.Item($r, 2).Value = $rowValues[1]
    $ws.Cells.Item($r, 3).Value = $rowValues[2]
}
